# Update Madigan bike hours (Riders / Average) on the Ridership sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Row 2 (Monday)
$ws.Range("C2").Value = 208
$ws.Range("D2").Value = 229.2

# Row 3 (Tuesday)
$ws.Range("C3").Value = 242
$ws.Range("D3").Value = 208.5

# Row 4 (Wednesday)
$ws.Range("C4").Value = 264
$ws.Range("D4").Value = 203.62

# Row 5 (Thursday)
$ws.Range("C5").Value = 282
$ws.Range("D5").Value = 229.67

# Row 6 (Friday)
$ws.Range("C6").Value = 246
$ws.Range("D6").Value = 238.69

# Row 7 (Saturday)
$ws.Range("C7").Value = 143
$ws.Range("D7").Value = 120.06

# Row 8 (Sunday)
$ws.Range("C8").Value = 96
$ws.Range("D8").Value = 102.47

$wb.Save()
